$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 15,15
$arr[0,0] = 9.020586599984323
$arr[0,1] = 24.67363666426199
$arr[0,2] = 10.72679972907323
$arr[0,3] = 15.01983100855659
$arr[0,4] = 5.363862443116282
$arr[0,5] = 15.55886854536925
$arr[0,6] = 6.775522177955993
$arr[0,7] = 15.23875041602943
$arr[0,8] = 3.848465881330531
$arr[0,9] = 9.698976488619106
$arr[0,10] = 10.1963451399622
$arr[0,11] = 26.8315523206924
$arr[0,12] = 12.1492242452136
$arr[0,13] = 8.712604807791573
$arr[0,14] = 7.617748177156455
$arr[1,0] = 6.65664162045891
$arr[1,1] = 26.46134906601895
$arr[1,2] = 10.36247134376251
$arr[1,3] = 13.720335414595
$arr[1,4] = 13.56880684472409
$arr[1,5] = 27.93170824760149
$arr[1,6] = 6.777624991253685
$arr[1,7] = 4.374621048645491
$arr[1,8] = 4.374621048645491
$arr[1,9] = 14.83870739332397
$arr[1,10] = 10.47634216671462
$arr[1,11] = 8.088321850007949
$arr[1,12] = 3.405432908434642
$arr[1,13] = 13.6890454762168
$arr[1,14] = 11.32278651377574
$arr[2,0] = 19.36468565236299
$arr[2,1] = 2.67685778436678
$arr[2,2] = 27.77094357738926
$arr[2,3] = 10.02187325172362
$arr[2,4] = 20.9879402421614
$arr[2,5] = 26.98507316143293
$arr[2,6] = 1.219657148093041
$arr[2,7] = 27.61381135219987
$arr[2,8] = 18.63592430642453
$arr[2,9] = 23.00091471856156
$arr[2,10] = 7.173679626033062
$arr[2,11] = 14.0096700652081
$arr[2,12] = 11.40699549279769
$arr[2,13] = 20.47118459153531
$arr[2,14] = 11.55551187777098
$arr[3,0] = 22.68498886286671
$arr[3,1] = 9.84480156645207
$arr[3,2] = 7.75567292966933
$arr[3,3] = 25.57184698124054
$arr[3,4] = 3.619691902740333
$arr[3,5] = 8.770812804112415
$arr[3,6] = 18.76423327876132
$arr[3,7] = 26.8533117801381
$arr[3,8] = 27.11544287115327
$arr[3,9] = 20.695752261988
$arr[3,10] = 12.41609532982015
$arr[3,11] = 12.41609532982015
$arr[3,12] = 19.79009900217976
$arr[3,13] = 12.08944564547694
$arr[3,14] = 10.65307612424885
$arr[4,0] = 9.035511301689411
$arr[4,1] = 9.035511301689411
$arr[4,2] = 4.642220944286654
$arr[4,3] = 4.642220944286654
$arr[4,4] = 13.21727771870139
$arr[4,5] = 24.20714636782846
$arr[4,6] = 27.66177531389596
$arr[4,7] = 24.49486028515228
$arr[4,8] = 22.9469062548071
$arr[4,9] = 7.242197118844668
$arr[4,10] = 4.942554006271341
$arr[4,11] = 4.942554006271341
$arr[4,12] = 5.360146629905638
$arr[4,13] = 25.29591658572435
$arr[4,14] = 15.90967438773843
$arr[5,0] = 7.763831849937518
$arr[5,1] = 17.56779647748941
$arr[5,2] = 17.56779647748941
$arr[5,3] = 13.99473355653367
$arr[5,4] = 8.775473471634051
$arr[5,5] = 13.89838365474218
$arr[5,6] = 10.9087375309873
$arr[5,7] = 13.11921595894084
$arr[5,8] = 23.37258051501598
$arr[5,9] = 3.640464006761429
$arr[5,10] = 8.853277735686179
$arr[5,11] = 20.56559500755781
$arr[5,12] = 15.3478589857668
$arr[5,13] = 22.95104756365309
$arr[5,14] = 20.9414362031996
$arr[6,0] = 27.23910782517492
$arr[6,1] = 22.953313223212
$arr[6,2] = 7.33513798986549
$arr[6,3] = 25.90703012061408
$arr[6,4] = 27.60899868009487
$arr[6,5] = 19.95831306889844
$arr[6,6] = 6.850920128876277
$arr[6,7] = 5.8312788998101
$arr[6,8] = 2.32234814199044
$arr[6,9] = 16.52020792562613
$arr[6,10] = 6.733267294665362
$arr[6,11] = 12.68573470767679
$arr[6,12] = 23.38219107985329
$arr[6,13] = 18.67833087102709
$arr[6,14] = 21.07392016333632
$arr[7,0] = 17.00461853819624
$arr[7,1] = 25.79513173688258
$arr[7,2] = 8.689696172901009
$arr[7,3] = 29.18581842389854
$arr[7,4] = 1.023208687771568
$arr[7,5] = 22.51138243461249
$arr[7,6] = 18.83166186437721
$arr[7,7] = 24.89218327056318
$arr[7,8] = 28.42345766063548
$arr[7,9] = 1.870231295868793
$arr[7,10] = 16.73684225641082
$arr[7,11] = 18.71471164207669
$arr[7,12] = 18.71471164207669
$arr[7,13] = 12.154105615584
$arr[7,14] = 13.93969278542829
$arr[8,0] = 15.06632607218783
$arr[8,1] = 20.50219241933116
$arr[8,2] = 27.08730749569193
$arr[8,3] = 27.17334500597135
$arr[8,4] = 11.17127773440506
$arr[8,5] = 26.55865803871687
$arr[8,6] = 6.027803756825253
$arr[8,7] = 24.60021765020461
$arr[8,8] = 16.5798847966854
$arr[8,9] = 12.33004800590377
$arr[8,10] = 9.021599877183981
$arr[8,11] = 1.503178212032449
$arr[8,12] = 29.59404721513443
$arr[8,13] = 13.12232191046616
$arr[8,14] = 21.77207320915707
$arr[9,0] = 22.50651870301571
$arr[9,1] = 9.054076096667702
$arr[9,2] = 17.79480919293188
$arr[9,3] = 5.219185239100367
$arr[9,4] = 19.87793269942794
$arr[9,5] = 19.54436703917215
$arr[9,6] = 5.827975690417865
$arr[9,7] = 7.576794605258355
$arr[9,8] = 5.161405084976556
$arr[9,9] = 5.884250129266718
$arr[9,10] = 10.04838488116379
$arr[9,11] = 13.71965238131274
$arr[9,12] = 13.71965238131274
$arr[9,13] = 9.395146774475188
$arr[9,14] = 18.35088276762268
$arr[10,0] = 14.77790663554525
$arr[10,1] = 15.18508406849698
$arr[10,2] = 7.646088667563495
$arr[10,3] = 9.074873529146771
$arr[10,4] = 18.44302727460981
$arr[10,5] = 2.725042060655934
$arr[10,6] = 16.75851609283128
$arr[10,7] = 14.71429577832866
$arr[10,8] = 2.606251532508013
$arr[10,9] = 21.79207171607265
$arr[10,10] = 6.812736040469777
$arr[10,11] = 6.812736040469777
$arr[10,12] = 8.334001432775583
$arr[10,13] = 28.10699257131816
$arr[10,14] = 15.67315567957748
$arr[11,0] = 3.602648557520952
$arr[11,1] = 1.127995291397784
$arr[11,2] = 4.194388733105692
$arr[11,3] = 25.16342597355398
$arr[11,4] = 26.62883943638965
$arr[11,5] = 16.31454813146462
$arr[11,6] = 21.77451527609803
$arr[11,7] = 10.70658040653859
$arr[11,8] = 27.23528548893918
$arr[11,9] = 20.28146964982595
$arr[11,10] = 27.08988582703597
$arr[11,11] = 28.08781533169034
$arr[11,12] = 19.6709025065967
$arr[11,13] = 18.1280946588224
$arr[11,14] = 29.29065843955219
$arr[12,0] = 18.74030233037773
$arr[12,1] = 25.05242669187031
$arr[12,2] = 4.689495516043783
$arr[12,3] = 5.647451961491337
$arr[12,4] = 4.556009600308807
$arr[12,5] = 4.580853601822001
$arr[12,6] = 13.65995247700404
$arr[12,7] = 22.2318788167087
$arr[12,8] = 21.92133394300387
$arr[12,9] = 29.6468141410543
$arr[12,10] = 1.862281673418898
$arr[12,11] = 7.450948809802339
$arr[12,12] = 6.18266254991766
$arr[12,13] = 6.18266254991766
$arr[12,14] = 20.29969975207337
$arr[13,0] = 5.528112844812471
$arr[13,1] = 2.918892749886233
$arr[13,2] = 15.35304944013477
$arr[13,3] = 15.96787024711226
$arr[13,4] = 4.990799186621143
$arr[13,5] = 1.976925898324472
$arr[13,6] = 29.81875034483644
$arr[13,7] = 17.75131314780145
$arr[13,8] = 1.360135267890718
$arr[13,9] = 21.10514726853357
$arr[13,10] = 6.203397730755318
$arr[13,11] = 6.407401983505927
$arr[13,12] = 3.988945583311863
$arr[13,13] = 25.58000586983143
$arr[13,14] = 17.74336195098626
$arr[14,0] = 29.54224791091988
$arr[14,1] = 9.350808393173704
$arr[14,2] = 21.8610039664438
$arr[14,3] = 11.49988175082895
$arr[14,4] = 25.17638713427582
$arr[14,5] = 24.8896393474005
$arr[14,6] = 17.66448751826199
$arr[14,7] = 22.17169911999229
$arr[14,8] = 11.44868701914386
$arr[14,9] = 11.3171204551866
$arr[14,10] = 26.65498224445452
$arr[14,11] = 24.67205453325809
$arr[14,12] = 9.796212077255566
$arr[14,13] = 9.865745869666636
$arr[14,14] = 8.019449013036297

$ws.Range("A1:O15").Value = $arr
